$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 35047.152
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 36429.04
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 109287.12
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -109623.12
$ws.Range("H40").Value = 4241.5835
$ws.Range("I40").Value = 3183.3333
$ws.Range("J40").Value = 5299.8335
$ws.Range("K40").Value = 3183.3333
$ws.Range("L40").Value = 5299.8335
$ws.Range("M40").Value = -3008.3333
$ws.Range("N40").Value = -5649.8335
$ws.Range("H96").Value = 94260
$ws.Range("I96").Value = 169287.5
$ws.Range("J96").Value = 4227
$ws.Range("K96").Value = 507862.5
$ws.Range("L96").Value = 12681
$ws.Range("M96").Value = -506489.5
$ws.Range("N96").Value = -15427
$ws.Range("H107").Value = 59781.383
$ws.Range("I107").Value = 294.83334
$ws.Range("K107").Value = 294.83334
$ws.Range("M107").Value = 1625.16666
$ws.Range("H132").Value = 2833.6943
$ws.Range("J132").Value = 30277.5
$ws.Range("L132").Value = 90832.5
$ws.Range("N132").Value = -95892.5
$ws.Range("H137").Value = 3134.7354
$ws.Range("I137").Value = 1187.909
$ws.Range("J137").Value = 6703.9165
$ws.Range("K137").Value = 3563.727
$ws.Range("L137").Value = 20111.7495
$ws.Range("M137").Value = -1013.727
$ws.Range("N137").Value = -25211.7495
$ws.Range("H138").Value = 2679.3677
$ws.Range("I138").Value = 1376.1875
$ws.Range("J138").Value = 3080.3462
$ws.Range("K138").Value = 4128.5625
$ws.Range("L138").Value = 9241.0386
$ws.Range("M138").Value = 1011.4375
$ws.Range("N138").Value = -19521.0386
$ws.Range("H141").Value = 40600.043
$ws.Range("I141").Value = 44370.477
$ws.Range("J141").Value = 1010.5
$ws.Range("K141").Value = 133111.431
$ws.Range("L141").Value = 3031.5
$ws.Range("M141").Value = -127931.431
$ws.Range("N141").Value = -13391.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14119936
$ws.Range("I32").Value = 3407157
$ws.Range("K32").Value = 3407157
$ws.Range("M32").Value = -3406870
$ws.Range("H45").Value = 14789
$ws.Range("I45").Value = 14541.471
$ws.Range("K45").Value = 14541.471
$ws.Range("M45").Value = -14164.471
$ws.Range("H61").Value = 1787.1111
$ws.Range("I61").Value = 1820.04
$ws.Range("K61").Value = 1820.04
$ws.Range("M61").Value = -1608.04
$ws.Range("H74").Value = 1586.0588
$ws.Range("I74").Value = 1590.8148
$ws.Range("K74").Value = 1590.8148
$ws.Range("M74").Value = -716.8148000000001
$ws.Range("H77").Value = 1586.0588
$ws.Range("I77").Value = 1590.8148
$ws.Range("K77").Value = 7954.074000000001
$ws.Range("M77").Value = -3586.074000000001
$ws.Range("H110").Value = 1530.2941
$ws.Range("I110").Value = 963
$ws.Range("K110").Value = 963
$ws.Range("M110").Value = 1082
$ws.Range("H132").Value = 5450.4443
$ws.Range("I132").Value = 5465.125
$ws.Range("K132").Value = 16395.375
$ws.Range("M132").Value = -13865.375
$ws.Range("H136").Value = 1787.1111
$ws.Range("I136").Value = 1820.04
$ws.Range("K136").Value = 5460.12
$ws.Range("M136").Value = -2910.12
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4510.553
$ws.Range("I134").Value = 2832.8965
$ws.Range("J134").Value = 7213.4443
$ws.Range("K134").Value = 8498.6895
$ws.Range("L134").Value = 21640.3329
$ws.Range("M134").Value = -5963.6895
$ws.Range("N134").Value = -26710.3329
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 75000
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H58").Value = 957.1429000000001
$ws.Range("I58").Value = 966.6667
$ws.Range("J58").Value = 900
$ws.Range("K58").Value = 966.6667
$ws.Range("L58").Value = 900
$ws.Range("M58").Value = -763.6667
$ws.Range("N58").Value = -1306
$ws.Range("H132").Value = 2828.7036
$ws.Range("I132").Value = 2573.3809
$ws.Range("J132").Value = 3722.3333
$ws.Range("K132").Value = 7720.1427
$ws.Range("L132").Value = 11166.9999
$ws.Range("M132").Value = -5190.1427
$ws.Range("N132").Value = -16226.9999
$ws.Range("H136").Value = 957.1429000000001
$ws.Range("I136").Value = 966.6667
$ws.Range("J136").Value = 900
$ws.Range("K136").Value = 2900.0001
$ws.Range("L136").Value = 2700
$ws.Range("M136").Value = -350.0001000000002
$ws.Range("N136").Value = -7800
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 11219.5
$ws.Range("I133").Value = 5878
$ws.Range("K133").Value = 17634
$ws.Range("M133").Value = -12574
$ws.Range("H137").Value = 4341.6665
$ws.Range("I137").Value = 1976.25
$ws.Range("J137").Value = 5201.8184
$ws.Range("K137").Value = 5928.75
$ws.Range("L137").Value = 15605.4552
$ws.Range("M137").Value = -828.75
$ws.Range("N137").Value = -25805.4552
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 18891966
$ws.Range("I11").Value = 1718800
$ws.Range("J11").Value = 26906110
$ws.Range("K11").Value = 1718800
$ws.Range("L11").Value = 26906110
$ws.Range("M11").Value = -1718661
$ws.Range("N11").Value = -26906388
$ws.Range("H14").Value = 65000260
$ws.Range("I14").Value = 65000260
$ws.Range("K14").Value = 65000260
$ws.Range("M14").Value = -65000092
$ws.Range("H18").Value = 8000
$ws.Range("J18").Value = 8000
$ws.Range("L18").Value = 8000
$ws.Range("N18").Value = -8586
$ws.Range("H70").Value = 7941.76
$ws.Range("I70").Value = 6109.0625
$ws.Range("J70").Value = 11199.889
$ws.Range("K70").Value = 6109.0625
$ws.Range("L70").Value = 11199.889
$ws.Range("M70").Value = -5839.0625
$ws.Range("N70").Value = -11739.889
$ws.Range("H73").Value = 7941.76
$ws.Range("I73").Value = 6109.0625
$ws.Range("J73").Value = 11199.889
$ws.Range("K73").Value = 6109.0625
$ws.Range("L73").Value = 11199.889
$ws.Range("M73").Value = -5173.0625
$ws.Range("N73").Value = -13071.889
$ws.Range("H132").Value = 5071.3438
$ws.Range("I132").Value = 5131.577
$ws.Range("J132").Value = 4810.3335
$ws.Range("K132").Value = 15394.731
$ws.Range("L132").Value = 14431.0005
$ws.Range("M132").Value = -12864.731
$ws.Range("N132").Value = -19491.0005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 1500
$ws.Range("I25").Value = 1500
$ws.Range("K25").Value = 1500
$ws.Range("M25").Value = -1270
$ws.Range("H46").Value = 957.4375
$ws.Range("I46").Value = 1046.5
$ws.Range("J46").Value = 904
$ws.Range("K46").Value = 1046.5
$ws.Range("L46").Value = 904
$ws.Range("M46").Value = -858.5
$ws.Range("N46").Value = -1280
$ws.Range("H55").Value = 303.58975
$ws.Range("I55").Value = 248.88889
$ws.Range("J55").Value = 350.4762
$ws.Range("K55").Value = 248.88889
$ws.Range("L55").Value = 350.4762
$ws.Range("M55").Value = -75.88889
$ws.Range("N55").Value = -696.4762000000001
$ws.Range("H69").Value = 59999.668
$ws.Range("J69").Value = 59999.668
$ws.Range("L69").Value = 59999.668
$ws.Range("N69").Value = -61621.668
$ws.Range("H72").Value = 59999.668
$ws.Range("J72").Value = 59999.668
$ws.Range("L72").Value = 179999.004
$ws.Range("N72").Value = -188111.004
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 19147.273
$ws.Range("J70").Value = 19147.273
$ws.Range("L70").Value = 19147.273
$ws.Range("N70").Value = -19777.273
$ws.Range("H73").Value = 19147.273
$ws.Range("J73").Value = 19147.273
$ws.Range("L73").Value = 19147.273
$ws.Range("N73").Value = -21331.273
$ws.Range("H81").Value = 83334660
$ws.Range("I81").Value = 90910350
$ws.Range("K81").Value = 181820700
$ws.Range("M81").Value = -181819639
$ws.Range("H84").Value = 83334660
$ws.Range("I84").Value = 90910350
$ws.Range("K84").Value = 909103500
$ws.Range("M84").Value = -909098196
$ws.Range("H107").Value = 26321546
$ws.Range("I107").Value = 8068.385
$ws.Range("J107").Value = 83334080
$ws.Range("K107").Value = 24205.155
$ws.Range("L107").Value = 250002240
$ws.Range("M107").Value = -22285.155
$ws.Range("N107").Value = -250006080
$ws.Range("H132").Value = 7955.143
$ws.Range("I132").Value = 8716.4
$ws.Range("J132").Value = 6052
$ws.Range("K132").Value = 26149.2
$ws.Range("L132").Value = 18156
$ws.Range("M132").Value = -23619.2
$ws.Range("N132").Value = -23216
